$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New "Deprivation" section (entered in the same order as the original
# author so the shared-string table lines up: the two measure
# descriptions first, then the section header).
$ws.Range("A16").Value = "% living in income deprived households reliant on means tested benefit"
$ws.Range("A16").WrapText = $true
$ws.Rows.Item(16).RowHeight = 30

$ws.Range("A17").Value = "% of people aged over 60 who live in pension credit households"
$ws.Range("A17").WrapText = $true
$ws.Rows.Item(17).RowHeight = 30

$ws.Range("A15").Value = "Deprivation - 2 columns:"

# New "English speaking abilities" row
$ws.Range("A19").Value = "English speaking abilities"
$ws.Range("B19").Value = "4 measure of language skills in households as requested by Rishi"

# Move the active selection like the source workbook
$ws.Range("A21").Select()
